$d = $word.ActiveDocument
$d.Content.Find.Execute("481×7=", $true, $false, $false, $false, $false, $true, 1, $false, "272×9=", 2) | Out-Null
$d.Content.Find.Execute("982×9=", $true, $false, $false, $false, $false, $true, 1, $false, "199×2=", 2) | Out-Null
$d.Content.Find.Execute("241×4=", $true, $false, $false, $false, $false, $true, 1, $false, "165×3=", 2) | Out-Null
$d.Content.Find.Execute("959×3=", $true, $false, $false, $false, $false, $true, 1, $false, "244×2=", 2) | Out-Null
$d.Content.Find.Execute("428×9=", $true, $false, $false, $false, $false, $true, 1, $false, "732×3=", 2) | Out-Null
$d.Content.Find.Execute("368×8=", $true, $false, $false, $false, $false, $true, 1, $false, "636×8=", 2) | Out-Null
$d.Content.Find.Execute("548×7=", $true, $false, $false, $false, $false, $true, 1, $false, "356×9=", 2) | Out-Null
$d.Content.Find.Execute("738×6=", $true, $false, $false, $false, $false, $true, 1, $false, "566×7=", 2) | Out-Null
$d.Content.Find.Execute("644×9=", $true, $false, $false, $false, $false, $true, 1, $false, "326×9=", 2) | Out-Null
$d.Content.Find.Execute("505×5=", $true, $false, $false, $false, $false, $true, 1, $false, "149×3=", 2) | Out-Null
$d.Content.Find.Execute("530×2=", $true, $false, $false, $false, $false, $true, 1, $false, "332×4=", 2) | Out-Null
$d.Content.Find.Execute("585×8=", $true, $false, $false, $false, $false, $true, 1, $false, "768×5=", 2) | Out-Null
$d.Content.Find.Execute("313×9=", $true, $false, $false, $false, $false, $true, 1, $false, "493×7=", 2) | Out-Null
$d.Content.Find.Execute("198×2=", $true, $false, $false, $false, $false, $true, 1, $false, "114×9=", 2) | Out-Null
$d.Content.Find.Execute("989×6=", $true, $false, $false, $false, $false, $true, 1, $false, "804×4=", 2) | Out-Null
$d.Content.Find.Execute("240×3=", $true, $false, $false, $false, $false, $true, 1, $false, "293×8=", 2) | Out-Null
$d.Content.Find.Execute("386×2=", $true, $false, $false, $false, $false, $true, 1, $false, "547×3=", 2) | Out-Null
$d.Content.Find.Execute("775×4=", $true, $false, $false, $false, $false, $true, 1, $false, "194×7=", 2) | Out-Null
$d.Content.Find.Execute("838×7=", $true, $false, $false, $false, $false, $true, 1, $false, "329×2=", 2) | Out-Null
$d.Content.Find.Execute("549×3=", $true, $false, $false, $false, $false, $true, 1, $false, "615×2=", 2) | Out-Null
$d.Content.Find.Execute("678×6=", $true, $false, $false, $false, $false, $true, 1, $false, "430×7=", 2) | Out-Null
$d.Content.Find.Execute("431×6=", $true, $false, $false, $false, $false, $true, 1, $false, "808×5=", 2) | Out-Null
$d.Content.Find.Execute("903×8=", $true, $false, $false, $false, $false, $true, 1, $false, "852×3=", 2) | Out-Null
$d.Content.Find.Execute("569×7=", $true, $false, $false, $false, $false, $true, 1, $false, "979×2=", 2) | Out-Null
$d.Content.Find.Execute("588×2=", $true, $false, $false, $false, $false, $true, 1, $false, "671×3=", 2) | Out-Null
